$d = $word.ActiveDocument

# The footer of the page used to contain, right after the
# "LOM3013: Ciencia dos Materiais (Requisito fraco)" requirement line:
#   - a blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# These three paragraphs (the site-navigation/footer boilerplate) were
# removed from the published page, so delete them here too, leaving the
# requirement line followed directly by the trailing blank paragraph and
# the page-break paragraph.

$paras = $d.Paragraphs
$count = $paras.Count

$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $text = $paras.Item($i).Range.Text
    if ($text -like "*LOM3013*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $nextText = $paras.Item($targetIndex + 2).Range.Text
    $afterText = $paras.Item($targetIndex + 3).Range.Text

    if (($nextText -like "*Ver no Jupiter*") -and ($afterText -like "*2020*")) {
        $startPara = $paras.Item($targetIndex + 1)
        $endPara = $paras.Item($targetIndex + 3)

        $start = $startPara.Range.Start
        $end = $endPara.Range.End

        $range = $d.Range($start, $end)
        $range.Delete()
    }
}
